$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = 12
$ws.Range("B13").Value = "teste"

$ws.Range("B14").Select()
